$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-23 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-24 Friday", 2) | Out-Null
$d.Content.Find.Execute("7+58=", $true, $false, $false, $false, $false, $true, 1, $false, "94-18=", 2) | Out-Null
$d.Content.Find.Execute("79-10=", $true, $false, $false, $false, $false, $true, 1, $false, "90-45=", 2) | Out-Null
$d.Content.Find.Execute("70-4=", $true, $false, $false, $false, $false, $true, 1, $false, "38-36=", 2) | Out-Null
$d.Content.Find.Execute("64+28=", $true, $false, $false, $false, $false, $true, 1, $false, "94-64=", 2) | Out-Null
$d.Content.Find.Execute("77-65=", $true, $false, $false, $false, $false, $true, 1, $false, "75-31=", 2) | Out-Null
$d.Content.Find.Execute("70-51=", $true, $false, $false, $false, $false, $true, 1, $false, "91+1=", 2) | Out-Null
$d.Content.Find.Execute("15+15=", $true, $false, $false, $false, $false, $true, 1, $false, "1+19=", 2) | Out-Null
$d.Content.Find.Execute("77-34=", $true, $false, $false, $false, $false, $true, 1, $false, "54-43=", 2) | Out-Null
$d.Content.Find.Execute("2+9=", $true, $false, $false, $false, $false, $true, 1, $false, "26-24=", 2) | Out-Null
$d.Content.Find.Execute("99-59=", $true, $false, $false, $false, $false, $true, 1, $false, "86-6=", 2) | Out-Null
$d.Content.Find.Execute("28+58=", $true, $false, $false, $false, $false, $true, 1, $false, "11+48=", 2) | Out-Null
$d.Content.Find.Execute("11-2=", $true, $false, $false, $false, $false, $true, 1, $false, "42+16=", 2) | Out-Null
$d.Content.Find.Execute("90-46=", $true, $false, $false, $false, $false, $true, 1, $false, "26-26=", 2) | Out-Null
$d.Content.Find.Execute("41-37=", $true, $false, $false, $false, $false, $true, 1, $false, "0+23=", 2) | Out-Null
$d.Content.Find.Execute("6+76=", $true, $false, $false, $false, $false, $true, 1, $false, "53-15=", 2) | Out-Null
$d.Content.Find.Execute("53+6=", $true, $false, $false, $false, $false, $true, 1, $false, "0+45=", 2) | Out-Null
$d.Content.Find.Execute("64+10=", $true, $false, $false, $false, $false, $true, 1, $false, "98-12=", 2) | Out-Null
$d.Content.Find.Execute("43+44=", $true, $false, $false, $false, $false, $true, 1, $false, "11+43=", 2) | Out-Null
$d.Content.Find.Execute("70-19=", $true, $false, $false, $false, $false, $true, 1, $false, "22+7=", 2) | Out-Null
$d.Content.Find.Execute("95-13=", $true, $false, $false, $false, $false, $true, 1, $false, "96+2=", 2) | Out-Null
$d.Content.Find.Execute("77+9=", $true, $false, $false, $false, $false, $true, 1, $false, "5+58=", 2) | Out-Null
$d.Content.Find.Execute("62-5=", $true, $false, $false, $false, $false, $true, 1, $false, "33+57=", 2) | Out-Null
$d.Content.Find.Execute("7+60=", $true, $false, $false, $false, $false, $true, 1, $false, "81-43=", 2) | Out-Null
$d.Content.Find.Execute("13+86=", $true, $false, $false, $false, $false, $true, 1, $false, "54-39=", 2) | Out-Null
$d.Content.Find.Execute("76+9=", $true, $false, $false, $false, $false, $true, 1, $false, "18+10=", 2) | Out-Null
$d.Content.Find.Execute("19+47=", $true, $false, $false, $false, $false, $true, 1, $false, "9+89=", 2) | Out-Null
$d.Content.Find.Execute("80-66=", $true, $false, $false, $false, $false, $true, 1, $false, "10+6=", 2) | Out-Null
$d.Content.Find.Execute("81-20=", $true, $false, $false, $false, $false, $true, 1, $false, "15+25=", 2) | Out-Null
$d.Content.Find.Execute("44+9=", $true, $false, $false, $false, $false, $true, 1, $false, "61-59=", 2) | Out-Null
$d.Content.Find.Execute("71+22=", $true, $false, $false, $false, $false, $true, 1, $false, "83-59=", 2) | Out-Null
$d.Content.Find.Execute("90-37=", $true, $false, $false, $false, $false, $true, 1, $false, "87-32=", 2) | Out-Null
$d.Content.Find.Execute("18+18=", $true, $false, $false, $false, $false, $true, 1, $false, "37+58=", 2) | Out-Null
$d.Content.Find.Execute("6+69=", $true, $false, $false, $false, $false, $true, 1, $false, "61+34=", 2) | Out-Null
$d.Content.Find.Execute("47+13=", $true, $false, $false, $false, $false, $true, 1, $false, "35+21=", 2) | Out-Null
$d.Content.Find.Execute("18+16=", $true, $false, $false, $false, $false, $true, 1, $false, "2+60=", 2) | Out-Null
$d.Content.Find.Execute("38+45=", $true, $false, $false, $false, $false, $true, 1, $false, "83-42=", 2) | Out-Null
$d.Content.Find.Execute("75+24=", $true, $false, $false, $false, $false, $true, 1, $false, "24-20=", 2) | Out-Null
$d.Content.Find.Execute("17+44=", $true, $false, $false, $false, $false, $true, 1, $false, "43-30=", 2) | Out-Null
$d.Content.Find.Execute("28+56=", $true, $false, $false, $false, $false, $true, 1, $false, "96-91=", 2) | Out-Null
$d.Content.Find.Execute("26+59=", $true, $false, $false, $false, $false, $true, 1, $false, "79+2=", 2) | Out-Null
$d.Content.Find.Execute("85-80=", $true, $false, $false, $false, $false, $true, 1, $false, "56+26=", 2) | Out-Null
$d.Content.Find.Execute("1+35=", $true, $false, $false, $false, $false, $true, 1, $false, "16-7=", 2) | Out-Null
$d.Content.Find.Execute("62-46=", $true, $false, $false, $false, $false, $true, 1, $false, "15+47=", 2) | Out-Null
$d.Content.Find.Execute("29+5=", $true, $false, $false, $false, $false, $true, 1, $false, "18+11=", 2) | Out-Null
$d.Content.Find.Execute("35+7=", $true, $false, $false, $false, $false, $true, 1, $false, "12+36=", 2) | Out-Null
$d.Content.Find.Execute("55-15=", $true, $false, $false, $false, $false, $true, 1, $false, "45+6=", 2) | Out-Null
$d.Content.Find.Execute("9-1=", $true, $false, $false, $false, $false, $true, 1, $false, "7+72=", 2) | Out-Null
$d.Content.Find.Execute("79-11=", $true, $false, $false, $false, $false, $true, 1, $false, "7+48=", 2) | Out-Null
$d.Content.Find.Execute("29-9=", $true, $false, $false, $false, $false, $true, 1, $false, "65-61=", 2) | Out-Null
$d.Content.Find.Execute("45-32=", $true, $false, $false, $false, $false, $true, 1, $false, "26+41=", 2) | Out-Null
$d.Content.Find.Execute("85-54=", $true, $false, $false, $false, $false, $true, 1, $false, "67-21=", 2) | Out-Null
$d.Content.Find.Execute("94-84=", $true, $false, $false, $false, $false, $true, 1, $false, "11+40=", 2) | Out-Null
$d.Content.Find.Execute("82-70=", $true, $false, $false, $false, $false, $true, 1, $false, "9+52=", 2) | Out-Null
$d.Content.Find.Execute("44-39=", $true, $false, $false, $false, $false, $true, 1, $false, "27-15=", 2) | Out-Null
$d.Content.Find.Execute("71-64=", $true, $false, $false, $false, $false, $true, 1, $false, "88-34=", 2) | Out-Null
$d.Content.Find.Execute("66+7=", $true, $false, $false, $false, $false, $true, 1, $false, "35+11=", 2) | Out-Null
$d.Content.Find.Execute("64-32=", $true, $false, $false, $false, $false, $true, 1, $false, "76-64=", 2) | Out-Null
$d.Content.Find.Execute("61+1=", $true, $false, $false, $false, $false, $true, 1, $false, "99-39=", 2) | Out-Null
$d.Content.Find.Execute("91-13=", $true, $false, $false, $false, $false, $true, 1, $false, "94-57=", 2) | Out-Null
$d.Content.Find.Execute("57+14=", $true, $false, $false, $false, $false, $true, 1, $false, "22-19=", 2) | Out-Null
$d.Content.Find.Execute("48+35=", $true, $false, $false, $false, $false, $true, 1, $false, "73-30=", 2) | Out-Null
$d.Content.Find.Execute("29+49=", $true, $false, $false, $false, $false, $true, 1, $false, "45+8=", 2) | Out-Null
$d.Content.Find.Execute("86-59=", $true, $false, $false, $false, $false, $true, 1, $false, "63+19=", 2) | Out-Null
$d.Content.Find.Execute("68-36=", $true, $false, $false, $false, $false, $true, 1, $false, "73+9=", 2) | Out-Null
$d.Content.Find.Execute("95-20=", $true, $false, $false, $false, $false, $true, 1, $false, "73-60=", 2) | Out-Null
$d.Content.Find.Execute("49+28=", $true, $false, $false, $false, $false, $true, 1, $false, "79-56=", 2) | Out-Null
$d.Content.Find.Execute("20+2=", $true, $false, $false, $false, $false, $true, 1, $false, "26+62=", 2) | Out-Null
$d.Content.Find.Execute("0+3=", $true, $false, $false, $false, $false, $true, 1, $false, "69-42=", 2) | Out-Null
$d.Content.Find.Execute("83-53=", $true, $false, $false, $false, $false, $true, 1, $false, "51+32=", 2) | Out-Null
$d.Content.Find.Execute("45+1=", $true, $false, $false, $false, $false, $true, 1, $false, "10+4=", 2) | Out-Null
$d.Content.Find.Execute("57+31=", $true, $false, $false, $false, $false, $true, 1, $false, "38+30=", 2) | Out-Null
$d.Content.Find.Execute("35-27=", $true, $false, $false, $false, $false, $true, 1, $false, "4+79=", 2) | Out-Null
$d.Content.Find.Execute("73-26=", $true, $false, $false, $false, $false, $true, 1, $false, "77-49=", 2) | Out-Null
$d.Content.Find.Execute("18+63=", $true, $false, $false, $false, $false, $true, 1, $false, "52-24=", 2) | Out-Null
$d.Content.Find.Execute("7-2=", $true, $false, $false, $false, $false, $true, 1, $false, "57-32=", 2) | Out-Null
$d.Content.Find.Execute("60+23=", $true, $false, $false, $false, $false, $true, 1, $false, "91-15=", 2) | Out-Null
$d.Content.Find.Execute("26+36=", $true, $false, $false, $false, $false, $true, 1, $false, "12+12=", 2) | Out-Null
$d.Content.Find.Execute("56-35=", $true, $false, $false, $false, $false, $true, 1, $false, "63-15=", 2) | Out-Null
$d.Content.Find.Execute("36+46=", $true, $false, $false, $false, $false, $true, 1, $false, "37-6=", 2) | Out-Null
$d.Content.Find.Execute("78-43=", $true, $false, $false, $false, $false, $true, 1, $false, "89-4=", 2) | Out-Null
$d.Content.Find.Execute("27-5=", $true, $false, $false, $false, $false, $true, 1, $false, "25-1=", 2) | Out-Null
$d.Content.Find.Execute("15+34=", $true, $false, $false, $false, $false, $true, 1, $false, "84-76=", 2) | Out-Null
$d.Content.Find.Execute("32+45=", $true, $false, $false, $false, $false, $true, 1, $false, "90+5=", 2) | Out-Null
$d.Content.Find.Execute("12+16=", $true, $false, $false, $false, $false, $true, 1, $false, "2+67=", 2) | Out-Null
$d.Content.Find.Execute("64-6=", $true, $false, $false, $false, $false, $true, 1, $false, "33+54=", 2) | Out-Null
$d.Content.Find.Execute("21+51=", $true, $false, $false, $false, $false, $true, 1, $false, "58-18=", 2) | Out-Null
$d.Content.Find.Execute("56-34=", $true, $false, $false, $false, $false, $true, 1, $false, "25-3=", 2) | Out-Null
$d.Content.Find.Execute("61-42=", $true, $false, $false, $false, $false, $true, 1, $false, "13+51=", 2) | Out-Null
$d.Content.Find.Execute("22+76=", $true, $false, $false, $false, $false, $true, 1, $false, "4+85=", 2) | Out-Null
$d.Content.Find.Execute("33+23=", $true, $false, $false, $false, $false, $true, 1, $false, "89-50=", 2) | Out-Null
$d.Content.Find.Execute("92-86=", $true, $false, $false, $false, $false, $true, 1, $false, "12+76=", 2) | Out-Null
$d.Content.Find.Execute("23+38=", $true, $false, $false, $false, $false, $true, 1, $false, "41+32=", 2) | Out-Null
$d.Content.Find.Execute("90-24=", $true, $false, $false, $false, $false, $true, 1, $false, "49-42=", 2) | Out-Null
$d.Content.Find.Execute("18-9=", $true, $false, $false, $false, $false, $true, 1, $false, "35+50=", 2) | Out-Null
$d.Content.Find.Execute("35+41=", $true, $false, $false, $false, $false, $true, 1, $false, "14+47=", 2) | Out-Null
$d.Content.Find.Execute("0+80=", $true, $false, $false, $false, $false, $true, 1, $false, "90-6=", 2) | Out-Null
$d.Content.Find.Execute("74-4=", $true, $false, $false, $false, $false, $true, 1, $false, "47+21=", 2) | Out-Null
$d.Content.Find.Execute("8+63=", $true, $false, $false, $false, $false, $true, 1, $false, "39+55=", 2) | Out-Null
$d.Content.Find.Execute("6-3=", $true, $false, $false, $false, $false, $true, 1, $false, "18-8=", 2) | Out-Null
$d.Content.Find.Execute("63+20=", $true, $false, $false, $false, $false, $true, 1, $false, "92-36=", 2) | Out-Null
